$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: remove indent from column B (data rows 2-38) ---
$ws.Range("B2:B38").IndentLevel = 0

# --- Formatting: add left horizontal alignment to column C (data rows 2-38) ---
$ws.Range("C2:C38").HorizontalAlignment = -4131
$ws.Range("B33").HorizontalAlignment = -4131

# --- New data row 39: EURAL code 200125 ---
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "200125"
$ws.Range("B39").Value = "spijsolie en -vetten"
$ws.Range("C39").Value = "Non-hazardous"
$ws.Range("A39:C39").HorizontalAlignment = -4131

# --- Update active selection to match author's last cursor position ---
$ws.Range("E17").Select()

Write-Host "done"
